$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 8. This shifts the old "totals" row (8)
#    down to row 9, and the old footer row (9) down to row 10 -- exactly
#    like Excel does when a new item row is added to the sale/shortage table.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Insert()

# Copy the formatting of the first item row (row 7) onto the freshly
# inserted row 8 so both item rows look identical.
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the row heights to match the final layout.
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 16.5

# Re-create the merges for the new item row 8 (same pattern as row 7).
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# The subtotal row (now row 9) keeps its P9:Q9 merge (re-applied defensively
# in case the shift above ever lost it).
$ws.Range("P9:Q9").Merge()

# ---------------------------------------------------------------------------
# 2. The item-name / counters columns are stored as text. Switch those
#    columns permanently to a text number format (numFmtId 49, "@"),
#    matching the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("Q7").NumberFormat = "@"

$ws.Range("C8:G8").NumberFormat = "@"
$ws.Range("N8:O8").NumberFormat = "@"
$ws.Range("H8:K8").NumberFormat = "@"
$ws.Range("Q8").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 3. Fill in the data for the two item rows.
# ---------------------------------------------------------------------------

# Row 7: FATROXIM 550 MG 30TAB
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "FATROXIM 550 MG 30TAB"
$ws.Range("H7").Value = "0:0"
$ws.Range("N7").Value = "360.00"
$ws.Range("Q7").Value = "0:1"

# L7 and P7 keep their original (numeric) number formats, but still hold
# text values, so flip the format to text only long enough to store the
# string, then restore the original numeric format code.
$fmtL7 = $ws.Range("L7").NumberFormat()
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "0"
$ws.Range("L7").NumberFormat = $fmtL7

$fmtP7 = $ws.Range("P7").NumberFormat()
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "118.8000"
$ws.Range("P7").NumberFormat = $fmtP7

# Row 8: NORHINOSE 50MCG/DOSE NASAL SPRAY 120 DOSES
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "NORHINOSE 50MCG/DOSE NASAL SPRAY 120 DOSES"
$ws.Range("H8").Value = "3:0"
$ws.Range("N8").Value = "90.00"
$ws.Range("Q8").Value = "1:0"

$fmtL8 = $ws.Range("L8").NumberFormat()
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "1"
$ws.Range("L8").NumberFormat = $fmtL8

$fmtP8 = $ws.Range("P8").NumberFormat()
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "90.0000"
$ws.Range("P8").NumberFormat = $fmtP8

# Row 9 (subtotal row): total of the two selling prices above.
$ws.Range("P9").Value = 208.80000000000001
$ws.Range("Q9").Value = ""

Write-Output "edit complete"
